# Insert one new daily price-observation row above the current row 81
# (Durazno / Florida King / Primera, 2020-12-04) so that every existing
# row from 81 down to 136 shifts down by one (new rows 82-137), and
# populate the freshly inserted row 81 with the new observation
# (Durazno / Early Majestic / Primera, 2021-11-12, Provincia de Limarí).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 81:136 down to 82:137, leaving a blank row 81 in place
# (formatting of the date column D is preserved from the row below).
$ws.Rows("81:81").Insert()

# Populate the new row 81 with the new data point.
$ws.Range("A81").Value = 10
$ws.Range("B81").Value = "Vega Modelo de Temuco"
$ws.Range("C81").Value = "La Araucanía"
$ws.Range("D81").Value = 44512
$ws.Range("E81").Value = 9
$ws.Range("F81").Value = "Fruta"
$ws.Range("G81").Value = 100103
$ws.Range("H81").Value = "Frutos de hueso (carozo)"
$ws.Range("I81").Value = 100103004
$ws.Range("J81").Value = "Durazno"
$ws.Range("K81").Value = "Early Majestic"
$ws.Range("L81").Value = "Primera"
$ws.Range("M81").Value = 65
$ws.Range("N81").Value = 34000
$ws.Range("O81").Value = 34000
$ws.Range("P81").Value = 34000
$ws.Range("Q81").Value = "$/bandeja 18 kilos granel"
$ws.Range("R81").Value = "Provincia de Limarí"
$ws.Range("S81").Value = 1889
$ws.Range("T81").Value = 18
